# Renamed con to ctr & added concentration symbol
#
# The "plot" column (C) used the abbreviation "con" for the control plots.
# Rename every occurrence of "con" to "ctr" so the shared-string table no
# longer needs the old label and gains the new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("timing")

# Find every cell in column C (the "plot" column) whose value is "con"
# and rename it to "ctr" - do this generically instead of hard-coding
# row numbers so it keeps working if the sheet's data range changes.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq "con") {
        $cell.Value = "ctr"
    }
}

# Mirror the author's cursor position captured in the saved workbook.
$ws.Range("C5").Select() | Out-Null
